# Add a new row (55) to Sheet1 with the match data for
# Sporting CP vs Rio Ave (2023-09-25), mirroring the formatting of the
# preceding data row (row 54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/borders/number formats) of the last existing
# data row down into the new row before writing values into it, so the new
# row matches the look of the rest of the table (bold/bordered index column,
# date-formatted match-date column, etc.).
$ws.Range("A54:V54").Copy()
$ws.Range("A55:V55").PasteSpecial(-4122)

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "portugal"
$ws.Range("C55").Value = "liga-portugal"
$ws.Range("D55").Value = "2023-2024"
$ws.Range("E55").Value = 45194.88541666666
$ws.Range("F55").Value = "Sporting CP"
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = "Rio Ave"
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1.18
$ws.Range("K55").Value = "18/09/2023 20:42"
$ws.Range("L55").Value = 1.17
$ws.Range("M55").Value = "25/09/2023 21:02"
$ws.Range("N55").Value = 7.68
$ws.Range("O55").Value = "18/09/2023 20:42"
$ws.Range("P55").Value = 7.84
$ws.Range("Q55").Value = "25/09/2023 21:02"
$ws.Range("R55").Value = 14.19
$ws.Range("S55").Value = "18/09/2023 20:42"
$ws.Range("T55").Value = 17.78
$ws.Range("U55").Value = "25/09/2023 21:02"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/sporting-lisbon-rio-ave/zkVaVVkM/"
